# Generate Report for Handback
# Adds a new handback row for file "f7324243-a653-4e8c-8e9e-f686619523f6.md"
# to the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$fileName   = "f7324243-a653-4e8c-8e9e-f686619523f6.md"
$baseId     = "f7324243-a653-4e8c-8e9e-f686619523f6"
$syncStatus = "Handed back: in sync with en-US"
$hlColor    = 15570276   # BGR COM color equivalent of RGB(0x64,0x95,0xED) -> matches existing FF6495ED hyperlink font

function Style-AsHyperlink($range) {
    # Re-apply the workbook's existing custom hyperlink font (underline, Calibri 11, #6495ED)
    # so the new cell's style dedupes against the pre-existing "HyperLink" cellXf instead of
    # Excel's own theme-colored built-in Hyperlink style.
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = $true
    $range.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (table: File Name / zh-cn / de-de) -> new row 5
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = $fileName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/95e6df379a01d58b025be581c1de74ce897c33d2/e2e/$fileName", "", "", $fileName) | Out-Null
Style-AsHyperlink $wsOverview.Range("A5")

$wsOverview.Range("B5").Value = $syncStatus
$wsOverview.Range("C5").Value = $syncStatus

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 5
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$baseId.b0dced80b65e25416406e3268a98ae1b620e7aa4.zh-cn.xlf"

$wsZhCn.Range("A5").Value = $fileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/95e6df379a01d58b025be581c1de74ce897c33d2/e2e/$fileName", "", "", $fileName) | Out-Null
Style-AsHyperlink $wsZhCn.Range("A5")

$wsZhCn.Range("B5").Value = $syncStatus

$wsZhCn.Range("C5").Value = $zhXlf
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f4428efa0e1a048676c8801b940034905590877/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null
Style-AsHyperlink $wsZhCn.Range("C5")

$wsZhCn.Range("D5").Value = "2016-03-10 04:49:14"
$wsZhCn.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("E5").Value = $fileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/40b737021eb206456a4de30004192601ee467abb/e2e/$fileName", "", "", $fileName) | Out-Null
Style-AsHyperlink $wsZhCn.Range("E5")

$wsZhCn.Range("F5").Value = $zhXlf
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/34df37d1655bf9cfcc4dbd7c8aed9e3ffed6c584/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null
Style-AsHyperlink $wsZhCn.Range("F5")

$wsZhCn.Range("G5").Value = "2016-03-10 04:49:52"
$wsZhCn.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H5").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 5
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$baseId.b0dced80b65e25416406e3268a98ae1b620e7aa4.de-de.xlf"

$wsDeDe.Range("A5").Value = $fileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/95e6df379a01d58b025be581c1de74ce897c33d2/e2e/$fileName", "", "", $fileName) | Out-Null
Style-AsHyperlink $wsDeDe.Range("A5")

$wsDeDe.Range("B5").Value = $syncStatus

$wsDeDe.Range("C5").Value = $deXlf
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0dc45b59b0ea48d1cc526dc38f69052a21ba3ad5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null
Style-AsHyperlink $wsDeDe.Range("C5")

$wsDeDe.Range("D5").Value = "2016-03-10 04:49:24"
$wsDeDe.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("E5").Value = $fileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ec3a0295bd752a6ff1dfa01004d3d71c2d390d8e/e2e/$fileName", "", "", $fileName) | Out-Null
Style-AsHyperlink $wsDeDe.Range("E5")

$wsDeDe.Range("F5").Value = $deXlf
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eb5365ecb6d13267108e54516137d1f48eebce78/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null
Style-AsHyperlink $wsDeDe.Range("F5")

$wsDeDe.Range("G5").Value = "2016-03-10 04:50:07"
$wsDeDe.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H5").Value = "Include"
